# AnalyticsCodeTables.xlsx edit script
# - Remove the Person sheet
# - Add more values to PersonRace (Hispanic) and reorder so Unknown is last
# - Add more values to BondType (Property, Deposit)
# - Make PersonAgeID (col A) the same values as AgeInYears (col B)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. PersonRace: insert "Hispanic" before the trailing "Unknown" entry,
#    so the final order is Asian, Black, American Indian, White, Hispanic, Unknown
# ---------------------------------------------------------------------
$raceWs = $wb.Worksheets.Item("PersonRace")
$raceWs.Range("A5").Value = 5
$raceWs.Range("B5").Value = "White"
$raceWs.Range("A6").Value = 6
$raceWs.Range("B6").Value = "Hispanic"
$raceWs.Range("A7").Value = 7
$raceWs.Range("B7").Value = "Unknown"

# ---------------------------------------------------------------------
# 2. BondType: add Property and Deposit as new bond types
# ---------------------------------------------------------------------
$bondWs = $wb.Worksheets.Item("BondType")
$bondWs.Range("A4").Value = 3
$bondWs.Range("B4").Value = "Property"
$bondWs.Range("A5").Value = 4
$bondWs.Range("B5").Value = "Deposit"

# ---------------------------------------------------------------------
# 3. PersonAge: make PersonAgeID (column A) match AgeInYears (column B)
# ---------------------------------------------------------------------
$ageWs = $wb.Worksheets.Item("PersonAge")
for ($row = 2; $row -le 34; $row++) {
    $ageWs.Cells.Item($row, 1).Value2 = $ageWs.Cells.Item($row, 2).Value2
}
$ageWs.Range("A35").Value = 51
$ageWs.Range("A36").Value = 52

# ---------------------------------------------------------------------
# 4. Remove the Person sheet (its data is now populated via demoBookingCount)
# ---------------------------------------------------------------------
$personWs = $wb.Worksheets.Item("Person")
$personWs.Delete()
